$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 18330051920122
$ws.Range("B2").Value = "GUZMAN"
$ws.Range("C2").Value = "MORALES"
$ws.Range("D2").Value = "OSWALDO IVAN"
$ws.Range("E2").Value = "DESARROLLA APLICACIONES MÓVILES PARA IOS"
$ws.Range("F2").Value = "6APV"
$ws.Range("G2").Value = 2
